$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Efnb2/Ephb1 sender-receiver table (rows 2-10) with the
# recomputed NATMI values, filling in the previously-missing FAPs target
# rows so every Sending x Target cluster combination (ECs/FAPs/sCs) is present.

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efnb2"
$ws.Cells.Item(2, 3).Value = "Ephb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 20.94432133333333
$ws.Cells.Item(2, 8).Value = 62.832964
$ws.Cells.Item(2, 9).Value = 0.7396577289668299
$ws.Cells.Item(2, 10).Value = 0.7396577289668298
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.8299853333333332
$ws.Cells.Item(2, 14).Value = 2.489956
$ws.Cells.Item(2, 15).Value = 0.7350786001848651
$ws.Cells.Item(2, 16).Value = 0.7350786001848651
$ws.Cells.Item(2, 17).Value = 17.38347952328711
$ws.Cells.Item(2, 18).Value = 156.451315709584
$ws.Cells.Item(2, 19).Value = 0.5437065680248537
$ws.Cells.Item(2, 20).Value = 0.5437065680248536

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efnb2"
$ws.Cells.Item(3, 3).Value = "Ephb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 20.94432133333333
$ws.Cells.Item(3, 8).Value = 62.832964
$ws.Cells.Item(3, 9).Value = 0.7396577289668299
$ws.Cells.Item(3, 10).Value = 0.7396577289668298
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.01102233333333333
$ws.Cells.Item(3, 14).Value = 0.033067
$ws.Cells.Item(3, 15).Value = 0.009761957268446888
$ws.Cells.Item(3, 16).Value = 0.009761957268446888
$ws.Cells.Item(3, 17).Value = 0.2308552911764445
$ws.Cells.Item(3, 18).Value = 2.077697620588
$ws.Cells.Item(3, 19).Value = 0.007220507143450664
$ws.Cells.Item(3, 20).Value = 0.007220507143450662

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efnb2"
$ws.Cells.Item(4, 3).Value = "Ephb1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 20.94432133333333
$ws.Cells.Item(4, 8).Value = 62.832964
$ws.Cells.Item(4, 9).Value = 0.7396577289668299
$ws.Cells.Item(4, 10).Value = 0.7396577289668298
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.2881033333333333
$ws.Cells.Item(4, 14).Value = 0.8643099999999999
$ws.Cells.Item(4, 15).Value = 0.2551594425466879
$ws.Cells.Item(4, 16).Value = 0.2551594425466879
$ws.Cells.Item(4, 17).Value = 6.034128790537778
$ws.Cells.Item(4, 18).Value = 54.30715911484
$ws.Cells.Item(4, 19).Value = 0.1887306537985255
$ws.Cells.Item(4, 20).Value = 0.1887306537985254

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efnb2"
$ws.Cells.Item(5, 3).Value = "Ephb1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.327094666666667
$ws.Cells.Item(5, 8).Value = 6.981284
$ws.Cells.Item(5, 9).Value = 0.08218235047311259
$ws.Cells.Item(5, 10).Value = 0.08218235047311258
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.8299853333333332
$ws.Cells.Item(5, 14).Value = 2.489956
$ws.Cells.Item(5, 15).Value = 0.7350786001848651
$ws.Cells.Item(5, 16).Value = 0.7350786001848651
$ws.Cells.Item(5, 17).Value = 1.931454442611556
$ws.Cells.Item(5, 18).Value = 17.383089983504
$ws.Cells.Item(5, 19).Value = 0.06041048714567759
$ws.Cells.Item(5, 20).Value = 0.06041048714567759

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efnb2"
$ws.Cells.Item(6, 3).Value = "Ephb1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.327094666666667
$ws.Cells.Item(6, 8).Value = 6.981284
$ws.Cells.Item(6, 9).Value = 0.08218235047311259
$ws.Cells.Item(6, 10).Value = 0.08218235047311258
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.01102233333333333
$ws.Cells.Item(6, 14).Value = 0.033067
$ws.Cells.Item(6, 15).Value = 0.009761957268446888
$ws.Cells.Item(6, 16).Value = 0.009761957268446888
$ws.Cells.Item(6, 17).Value = 0.02565001311422222
$ws.Cells.Item(6, 18).Value = 0.230850118028
$ws.Cells.Item(6, 19).Value = 0.000802260593539051
$ws.Cells.Item(6, 20).Value = 0.0008022605935390509

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efnb2"
$ws.Cells.Item(7, 3).Value = "Ephb1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.327094666666667
$ws.Cells.Item(7, 8).Value = 6.981284
$ws.Cells.Item(7, 9).Value = 0.08218235047311259
$ws.Cells.Item(7, 10).Value = 0.08218235047311258
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.2881033333333333
$ws.Cells.Item(7, 14).Value = 0.8643099999999999
$ws.Cells.Item(7, 15).Value = 0.2551594425466879
$ws.Cells.Item(7, 16).Value = 0.2551594425466879
$ws.Cells.Item(7, 17).Value = 0.6704437304488889
$ws.Cells.Item(7, 18).Value = 6.03399357404
$ws.Cells.Item(7, 19).Value = 0.02096960273389594
$ws.Cells.Item(7, 20).Value = 0.02096960273389594

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Efnb2"
$ws.Cells.Item(8, 3).Value = "Ephb1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 5.044818
$ws.Cells.Item(8, 8).Value = 15.134454
$ws.Cells.Item(8, 9).Value = 0.1781599205600575
$ws.Cells.Item(8, 10).Value = 0.1781599205600575
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.8299853333333332
$ws.Cells.Item(8, 14).Value = 2.489956
$ws.Cells.Item(8, 15).Value = 0.7350786001848651
$ws.Cells.Item(8, 16).Value = 0.7350786001848651
$ws.Cells.Item(8, 17).Value = 4.187124949336
$ws.Cells.Item(8, 18).Value = 37.684124544024
$ws.Cells.Item(8, 19).Value = 0.1309615450143339
$ws.Cells.Item(8, 20).Value = 0.1309615450143339

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Efnb2"
$ws.Cells.Item(9, 3).Value = "Ephb1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 5.044818
$ws.Cells.Item(9, 8).Value = 15.134454
$ws.Cells.Item(9, 9).Value = 0.1781599205600575
$ws.Cells.Item(9, 10).Value = 0.1781599205600575
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.01102233333333333
$ws.Cells.Item(9, 14).Value = 0.033067
$ws.Cells.Item(9, 15).Value = 0.009761957268446888
$ws.Cells.Item(9, 16).Value = 0.009761957268446888
$ws.Cells.Item(9, 17).Value = 0.055605665602
$ws.Cells.Item(9, 18).Value = 0.500450990418
$ws.Cells.Item(9, 19).Value = 0.001739189531457174
$ws.Cells.Item(9, 20).Value = 0.001739189531457174

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Efnb2"
$ws.Cells.Item(10, 3).Value = "Ephb1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 5.044818
$ws.Cells.Item(10, 8).Value = 15.134454
$ws.Cells.Item(10, 9).Value = 0.1781599205600575
$ws.Cells.Item(10, 10).Value = 0.1781599205600575
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.2881033333333333
$ws.Cells.Item(10, 14).Value = 0.8643099999999999
$ws.Cells.Item(10, 15).Value = 0.2551594425466879
$ws.Cells.Item(10, 16).Value = 0.2551594425466879
$ws.Cells.Item(10, 17).Value = 1.45342888186
$ws.Cells.Item(10, 18).Value = 13.08085993674
$ws.Cells.Item(10, 19).Value = 0.04545918601426648
$ws.Cells.Item(10, 20).Value = 0.04545918601426648
